$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "27.978.80"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").Value = "1.862.21"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  -0.37%  "
Set-TextCell "D5" "335.81"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("E6").Value = "  -0.39%  "
Set-TextCell "D7" "0.4701"
$ws.Range("E7").Value = "  +1.18%  "
Set-TextCell "D8" "0.3891"
$ws.Range("E8").Value = "  +0.93%  "
Set-TextCell "D9" "46.74"
$ws.Range("E9").Value = "  +1.32%  "
Set-TextCell "D10" "0.07964"
$ws.Range("E10").Value = "  +0.71%  "
Set-TextCell "D11" "0.9785"
$ws.Range("E11").Value = "  -2.03%  "
Set-TextCell "D12" "21.50"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "1.853.23"
$ws.Range("E13").Value = "  +0.08%  "
Set-TextCell "D14" "5.921"
$ws.Range("E14").Value = "  -0.82%  "
Set-TextCell "D15" "7.198"
Set-TextCell "D16" "91.52"
$ws.Range("E16").Value = "  +3.35%  "
Set-TextCell "D17" "1.003"
$ws.Range("E17").Value = "  -0.47%  "
Set-TextCell "D18" "0.00001039"
$ws.Range("E18").Value = "  +0.31%  "
Set-TextCell "D19" "0.06601"
$ws.Range("E19").Value = "  -1.21%  "
Set-TextCell "D20" "17.49"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "27.961.69"
$ws.Range("E22").Value = "  +1.39%  "
Set-TextCell "D23" "5.390"
$ws.Range("E23").Value = "  -0.27%  "
Set-TextCell "D24" "10.91"
$ws.Range("E24").Value = "  +0.26%  "
Set-TextCell "D25" "2.290"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").Value = "2.074.10"
$ws.Range("E26").Value = "  +0.15%  "
Set-TextCell "D27" "159.24"
$ws.Range("E27").Value = "  +0.22%  "
Set-TextCell "D28" "19.53"
$ws.Range("E28").Value = "  -0.01%  "
Set-TextCell "D29" "2.095"
$ws.Range("E29").Value = "  -1.05%  "
Set-TextCell "D30" "5.441"
$ws.Range("E30").Value = "  +0.44%  "
Set-TextCell "D31" "119.19"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D32" "0.9574"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D33" "0.09472"
$ws.Range("E33").Value = "  +0.68%  "
Set-TextCell "D34" "3.580"
$ws.Range("E34").Value = "  -0.49%  "
Set-TextCell "D35" "5.301"
$ws.Range("E35").Value = "  -0.12%  "
Set-TextCell "D36" "1.347"
$ws.Range("E36").Value = "  +0.40%  "
Set-TextCell "D37" "0.06085"
$ws.Range("E37").Value = "  +0.57%  "
Set-TextCell "D38" "0.02241"
$ws.Range("E38").Value = "  +0.38%  "
Set-TextCell "D39" "8.279"
$ws.Range("E39").Value = "  -0.41%  "
Set-TextCell "D40" "1.163"
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("E41").Value = "  -0.19%  "
Set-TextCell "D42" "0.5887"
$ws.Range("E42").Value = "  -0.31%  "
Set-TextCell "D43" "0.1863"
$ws.Range("E43").Value = "  -0.15%  "
Set-TextCell "D44" "10.19"
$ws.Range("E44").Value = "  -1.33%  "
Set-TextCell "D45" "1.282"
$ws.Range("E45").Value = "  +3.50%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell "D46" "0.5524"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D47" "12.11"
$ws.Range("E47").Value = "  -0.98%  "
Set-TextCell "D48" "1.943"
$ws.Range("E48").Value = "  +1.46%  "
Set-TextCell "D49" "0.06865"
$ws.Range("E49").Value = "  +2.40%  "
Set-TextCell "D50" "111.53"
$ws.Range("E50").Value = "  +0.88%  "
Set-TextCell "D51" "1.002"
$ws.Range("E51").Value = "  -32.83%  "
